# feat: add theme customization options for document generation
#
# Restyle the built-in Heading 1-6 styles with a new "theme" color
# palette, bold headings with explicit spacing/outline levels, tidy up
# the Hyperlink character style (drop the hard-coded underline color)
# and recolor the MdCode block style.

$d = $word.ActiveDocument

function Set-HeadingStyle {
    param([string]$Name, [int]$Color, [int]$Size, [bool]$Bold, [bool]$Italic, [double]$SpaceBeforePt, [double]$SpaceAfterPt, [int]$OutlineLevel)

    $s = $d.Styles.Item($Name)

    # Paragraph-level formatting: keep with next, before/after spacing,
    # outline level (used for the new TOC-ready heading levels).
    $s.ParagraphFormat.KeepWithNext = $true
    $s.ParagraphFormat.SpaceBefore = $SpaceBeforePt
    $s.ParagraphFormat.SpaceAfter = $SpaceAfterPt
    $s.ParagraphFormat.OutlineLevel = $OutlineLevel

    # Run-level formatting: new theme color + size (ascii/complex-script
    # kept in sync) + bold per the new theme spec. Italic is only ever
    # turned ON here -- styles that must lose italic are cleared by the
    # caller beforehand so we don't stamp a spurious explicit "off" flag
    # on styles that never had italic set in the first place.
    $s.Font.Bold = $Bold
    $s.Font.BoldBi = $Bold
    if ($Italic) {
        $s.Font.Italic = $true
        $s.Font.ItalicBi = $true
    }
    $s.Font.Color = $Color
    $s.Font.Size = $Size
    $s.Font.SizeBi = $Size
}

# Heading 1: #2F5597, 18pt, bold, keepNext, before 24pt / after 12pt, outline level 1
Set-HeadingStyle "Heading 1" 0x97552F 18 $true $false 24 12 1

# Heading 2: #5B9BD5, 16pt, bold, keepNext, before 20pt / after 10pt, outline level 2
Set-HeadingStyle "Heading 2" 0xD59B5B 16 $true $false 20 10 2

# Heading 3: #44546A, 14pt, bold, keepNext, before 16pt / after 8pt, outline level 3
Set-HeadingStyle "Heading 3" 0x6A5444 14 $true $false 16 8 3

# Heading 4: #44546A, 13pt, bold, italic removed, keepNext, before 14pt / after 7pt, outline level 4
$h4 = $d.Styles.Item("Heading 4")
$h4.Font.Italic = $false
$h4.Font.ItalicBi = $false
Set-HeadingStyle "Heading 4" 0x6A5444 13 $true $false 14 7 4

# Heading 5: #44546A, 12pt, bold + italic, keepNext, before 12pt / after 6pt, outline level 5
Set-HeadingStyle "Heading 5" 0x6A5444 12 $true $true 12 6 5

# Heading 6: #44546A, 12pt, explicit not-bold + italic, keepNext, before 12pt / after 6pt, outline level 6
Set-HeadingStyle "Heading 6" 0x6A5444 12 $false $true 12 6 6

# Hyperlink character style: drop the hard-coded underline color, let it
# follow the text color automatically.
$hyperlinkStyle = $d.Styles.Item("Hyperlink")
$hyperlinkStyle.Font.UnderlineColor = -16777216

# MdCode block style: recolor from green to the new crimson accent.
$mdCode = $d.Styles.Item("MdCode")
$mdCode.Font.Color = 0x4E25C7
